# Backup QR Scanner data - 10/05/2025, 3:16:16 AM
# Appends the newly scanned row (student 885) to the Biochemistry checklist.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophes force Excel to keep numeric/date-looking values
# ("885", "10/05/2025") stored as text, matching the rest of the sheet.
$ws.Range("A4").Value = "'885"
$ws.Range("B4").Value = "Biochemistry"
$ws.Range("C4").Value = "'10/05/2025"
$ws.Range("D4").Value = "03:16:00"
$ws.Range("E4").Value = "'885"

# Drop the quote-prefix formatting that the apostrophe trick leaves behind
# so the new row doesn't pick up an extra cell style vs. the existing rows.
$ws.Range("A4:E4").ClearFormats()
